# Auto-generated edits applied to Seraph_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for specific leve rows
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 176.53847
$ws.Range("I9").Value = 91.75
$ws.Range("J9").Value = 312.2
$ws.Range("K9").Value = 91.75
$ws.Range("L9").Value = 312.2
$ws.Range("M9").Value = 77.25
$ws.Range("N9").Value = -650.2

$ws.Range("H112").Value = 2200.3635
$ws.Range("J112").Value = 2262.7896
$ws.Range("L112").Value = 6788.3688
$ws.Range("N112").Value = -9004.3688

$ws.Range("H127").Value = 1379.8
$ws.Range("I127").Value = 1133
$ws.Range("J127").Value = 1750
$ws.Range("K127").Value = 3399
$ws.Range("L127").Value = 5250
$ws.Range("M127").Value = 1561
$ws.Range("N127").Value = -15170

$ws.Range("H129").Value = 2084.7856
$ws.Range("I129").Value = 798.55554
$ws.Range("K129").Value = 2395.66662
$ws.Range("M129").Value = 2604.33338

$ws.Range("H132").Value = 3836.5557
$ws.Range("I132").Value = 3790
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 11370
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -8840
$ws.Range("N132").Value = -17058.5

$ws.Range("H137").Value = 5872.95
$ws.Range("I137").Value = 6031.1333
$ws.Range("K137").Value = 18093.3999
$ws.Range("M137").Value = -15543.3999

$ws.Range("H138").Value = 8581.672
$ws.Range("I138").Value = 7869.2856
$ws.Range("J138").Value = 8781.139999999999
$ws.Range("K138").Value = 23607.8568
$ws.Range("L138").Value = 26343.42
$ws.Range("M138").Value = -18467.8568
$ws.Range("N138").Value = -36623.42

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19037.545
$ws.Range("I32").Value = 13193.591
$ws.Range("J32").Value = 30725.455
$ws.Range("K32").Value = 13193.591
$ws.Range("L32").Value = 30725.455
$ws.Range("M32").Value = -12906.591
$ws.Range("N32").Value = -31299.455

$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 2540.2354
$ws.Range("I132").Value = 1685.8667
$ws.Range("K132").Value = 5057.6001
$ws.Range("M132").Value = -2527.6001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4470.6665
$ws.Range("I86").Value = 4866.4443
$ws.Range("J86").Value = 3283.3333
$ws.Range("K86").Value = 4866.4443
$ws.Range("L86").Value = 3283.3333
$ws.Range("M86").Value = -3743.4443
$ws.Range("N86").Value = -5529.3333

$ws.Range("H89").Value = 4470.6665
$ws.Range("I89").Value = 4866.4443
$ws.Range("J89").Value = 3283.3333
$ws.Range("K89").Value = 24332.2215
$ws.Range("L89").Value = 16416.6665
$ws.Range("M89").Value = -18716.2215
$ws.Range("N89").Value = -27648.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 14705.267
$ws.Range("I99").Value = 15664.667
$ws.Range("J99").Value = 14465.417
$ws.Range("K99").Value = 15664.667
$ws.Range("L99").Value = 14465.417
$ws.Range("M99").Value = -14166.667
$ws.Range("N99").Value = -17461.417

$ws.Range("H105").Value = 5157.8887
$ws.Range("I105").Value = 3618
$ws.Range("K105").Value = 3618
$ws.Range("M105").Value = -1871

$ws.Range("H107").Value = 479.33334
$ws.Range("I107").Value = 310.83334
$ws.Range("J107").Value = 647.8333
$ws.Range("K107").Value = 310.83334
$ws.Range("L107").Value = 647.8333
$ws.Range("M107").Value = 1609.16666
$ws.Range("N107").Value = -4487.8333

$ws.Range("H109").Value = 64598.332
$ws.Range("J109").Value = 64598.332
$ws.Range("L109").Value = 64598.332
$ws.Range("N109").Value = -66678.33199999999

$ws.Range("H126").Value = 14705.267
$ws.Range("I126").Value = 15664.667
$ws.Range("J126").Value = 14465.417
$ws.Range("K126").Value = 46994.001
$ws.Range("L126").Value = 43396.251
$ws.Range("M126").Value = -44524.001
$ws.Range("N126").Value = -48336.251

$ws.Range("H132").Value = 2280.5945
$ws.Range("I132").Value = 2121.1
$ws.Range("J132").Value = 2964.1428
$ws.Range("K132").Value = 6363.299999999999
$ws.Range("L132").Value = 8892.428400000001
$ws.Range("M132").Value = -3833.299999999999
$ws.Range("N132").Value = -13952.4284

$ws.Range("H134").Value = 2635.6206
$ws.Range("I134").Value = 2166.6
$ws.Range("J134").Value = 3677.889
$ws.Range("K134").Value = 6499.799999999999
$ws.Range("L134").Value = 11033.667
$ws.Range("M134").Value = -3964.799999999999
$ws.Range("N134").Value = -16103.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 880
$ws.Range("I48").Value = 173.33333
$ws.Range("J48").Value = 3000
$ws.Range("K48").Value = 519.99999
$ws.Range("L48").Value = 9000
$ws.Range("M48").Value = -269.99999
$ws.Range("N48").Value = -9500

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9999
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 9999
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9999
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -11995

$ws.Range("H83").Value = 9999
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 9999
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 49995
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -59979

$ws.Range("H134").Value = 60630
$ws.Range("J134").Value = 60630
$ws.Range("L134").Value = 181890
$ws.Range("N134").Value = -186960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H55").Value = 3057
$ws.Range("I55").Value = 3057
$ws.Range("K55").Value = 3057
$ws.Range("M55").Value = -2884

$ws.Range("H122").Value = 9128.571
$ws.Range("I122").Value = 7966.6665
$ws.Range("K122").Value = 23899.9995
$ws.Range("M122").Value = -21449.9995

$ws.Range("H132").Value = 5116.5835
$ws.Range("I132").Value = 3224.75
$ws.Range("K132").Value = 9674.25
$ws.Range("M132").Value = -7144.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 65999.664
$ws.Range("J64").Value = 65999.664
$ws.Range("L64").Value = 65999.664
$ws.Range("N64").Value = -66495.664

$ws.Range("H67").Value = 65999.664
$ws.Range("J67").Value = 65999.664
$ws.Range("L67").Value = 65999.664
$ws.Range("N67").Value = -67715.664

$ws.Range("H81").Value = 3778.0715
$ws.Range("J81").Value = 7000
$ws.Range("L81").Value = 14000
$ws.Range("N81").Value = -16122

$ws.Range("H84").Value = 3778.0715
$ws.Range("J84").Value = 7000
$ws.Range("L84").Value = 70000
$ws.Range("N84").Value = -80608

$ws.Range("H107").Value = 1800.375
$ws.Range("I107").Value = 780.6
$ws.Range("K107").Value = 2341.8
$ws.Range("M107").Value = -421.8000000000002

$ws.Range("H132").Value = 1234.875
$ws.Range("I132").Value = 1044.75
$ws.Range("J132").Value = 1425
$ws.Range("K132").Value = 3134.25
$ws.Range("L132").Value = 4275
$ws.Range("M132").Value = -604.25
$ws.Range("N132").Value = -9335
